# Add account information rows when creating a new role; add init property module.
# Appends 41 new property rows (rows 12-52) to the "Property1" sheet, mirroring the
# formatting already used by row 11 (the last existing data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (relative to the "test<n>" property name) whose B column holds the literal
# text "1" (formatted like a group header, same as row 11) instead of the numeric
# value 1.
$stringBRows = @(16,21,26,31,33,38,43,45,50)

for ($r = 12; $r -le 52; $r++) {
    $n = $r - 10

    # Column A: property name "test<n>", formatted like A11.
    $ws.Range("A11").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = "test$n"

    if ($stringBRows -contains $r) {
        # Column B: literal text "1", formatted like B11.
        $ws.Range("B11").Copy()
        $ws.Range("B$r").PasteSpecial(-4122)
        $ws.Range("B$r").Value = "1"
    } else {
        # Column B: numeric 1, default column formatting.
        $ws.Range("B$r").Value = 1
    }
}

$ws.Range("B33").Select() | Out-Null
